# Removed Test Case Inter-Dependency
# Update the loan product name (productname) and short name (shortname)
# on the ProductLoanInput sheet so this workbook no longer depends on the
# identifiers used by another test case, and mirror the product name
# change on the ProductLoanOutput sheet (which carries its own copy of
# the same text, not a formula reference).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "2597-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-AMT-VAR-INST-FIX-INST-AMT-1st"
$newShortName = "259q"

# productname (row 1) on both sheets
$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

# shortname (row 2) on the input sheet - switches from a numeric literal
# to a text value, keeping its existing number-format style
$ws1.Range("B2").Value = $newShortName

# Move the active selection on the input sheet to B8
$ws1.Activate()
$ws1.Range("B8").Select()
